# Refresh the cryptos list with the latest scrape snapshot.
# Only the "Price" (D) and "Volume(1h)" (E) columns change; row order,
# coin names and links stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of "Price" cells are plain decimals (e.g. "154.14") that Excel
# would otherwise auto-detect as numbers on assignment. Force those cells
# to Text first so the updated price keeps being stored as plain text,
# exactly like the other Price cells (e.g. "60.326.96") that can only ever
# be text.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "60.326.96"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "2.592.79"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "510.12"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").Value = "154.14"
$ws.Range("E6").Value = "  -3.39%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -3.90%  "
$ws.Range("D9").Value = "2.599.42"
$ws.Range("E9").Value = "  -3.33%  "
$ws.Range("D10").Value = "6.70"
$ws.Range("E10").Value = "  +8.75%  "
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").Value = "3.045.85"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").Value = "60.330.49"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").Value = "21.59"
$ws.Range("E16").Value = "  -3.00%  "
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "2.594.83"
$ws.Range("E18").Value = "  -3.35%  "
$ws.Range("D19").Value = "4.75"
$ws.Range("E19").Value = "  -1.61%  "
$ws.Range("D20").Value = "352.21"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "10.55"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "60.23"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").Value = "0.421"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "0.0₃0842"
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  -3.18%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "19.39"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").Value = "152.12"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").Value = "5.73"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "4.01"
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("D37").Value = "0.860"
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("D38").Value = "1.48"
$ws.Range("E38").Value = "  -4.24%  "
$ws.Range("D39").Value = "0.844"
$ws.Range("E39").Value = "  -4.87%  "
$ws.Range("D40").Value = "36.08"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").Value = "3.75"
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("D42").Value = "300.33"
$ws.Range("E42").Value = "  -3.12%  "
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("E44").Value = "  -4.66%  "
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").Value = "0.0553"
$ws.Range("E46").Value = "  -4.48%  "
$ws.Range("E47").Value = "  -3.35%  "
$ws.Range("D48").Value = "4.84"
$ws.Range("E48").Value = "  -4.14%  "
$ws.Range("D49").Value = "0.0233"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").Value = "1.995.51"
$ws.Range("E51").Value = "  -2.25%  "
